$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 20585.95
$ws.Range("I98").Value = 820.4074000000001
$ws.Range("J98").Value = 61637.46
$ws.Range("K98").Value = 820.4074000000001
$ws.Range("L98").Value = 61637.46
$ws.Range("M98").Value = 677.5925999999999
$ws.Range("N98").Value = -64633.46
$ws.Range("H117").Value = 48510
$ws.Range("J117").Value = 48510
$ws.Range("L117").Value = 48510
$ws.Range("N117").Value = -57688
$ws.Range("H122").Value = 20585.95
$ws.Range("I122").Value = 820.4074000000001
$ws.Range("J122").Value = 61637.46
$ws.Range("K122").Value = 2461.2222
$ws.Range("L122").Value = 184912.38
$ws.Range("M122").Value = -11.22220000000016
$ws.Range("N122").Value = -189812.38
$ws.Range("H132").Value = 18623.908
$ws.Range("I132").Value = 2656.8262
$ws.Range("J132").Value = 110434.625
$ws.Range("K132").Value = 7970.4786
$ws.Range("L132").Value = 331303.875
$ws.Range("M132").Value = -5440.4786
$ws.Range("N132").Value = -336363.875

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H6").Value = 66250.5
$ws.Range("I6").Value = 40000.668
$ws.Range("J6").Value = 82000.39999999999
$ws.Range("K6").Value = 40000.668
$ws.Range("L6").Value = 82000.39999999999
$ws.Range("M6").Value = -39827.668
$ws.Range("N6").Value = -82346.39999999999
$ws.Range("H80").Value = 51394.145
$ws.Range("J80").Value = 51394.145
$ws.Range("L80").Value = 51394.145
$ws.Range("N80").Value = -53390.145
$ws.Range("H83").Value = 51394.145
$ws.Range("J83").Value = 51394.145
$ws.Range("L83").Value = 154182.435
$ws.Range("N83").Value = -164166.435
$ws.Range("H98").Value = 43355
$ws.Range("J98").Value = 43355
$ws.Range("L98").Value = 43355
$ws.Range("N98").Value = -49345
$ws.Range("H101").Value = 48598
$ws.Range("J101").Value = 48598
$ws.Range("L101").Value = 48598
$ws.Range("N101").Value = -55088
$ws.Range("H104").Value = 41441.668
$ws.Range("J104").Value = 41441.668
$ws.Range("L104").Value = 41441.668
$ws.Range("N104").Value = -48429.668
$ws.Range("H105").Value = 47942.4
$ws.Range("J105").Value = 47942.4
$ws.Range("L105").Value = 47942.4
$ws.Range("N105").Value = -54930.4
$ws.Range("H106").Value = 47090.5
$ws.Range("J106").Value = 47090.5
$ws.Range("L106").Value = 47090.5
$ws.Range("N106").Value = -49614.5
$ws.Range("H107").Value = 34829.332
$ws.Range("J107").Value = 34829.332
$ws.Range("L107").Value = 34829.332
$ws.Range("N107").Value = -42509.332
$ws.Range("H109").Value = 42512.8
$ws.Range("J109").Value = 42512.8
$ws.Range("L109").Value = 42512.8
$ws.Range("N109").Value = -45286.8
$ws.Range("H113").Value = 0
$ws.Range("J113").Value = 0
$ws.Range("N113").ClearContents()
$ws.Range("H114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("N114").ClearContents()
$ws.Range("H117").Value = 47895.832
$ws.Range("J117").Value = 47895.832
$ws.Range("L117").Value = 47895.832
$ws.Range("N117").Value = -57073.832
$ws.Range("H118").Value = 49344
$ws.Range("J118").Value = 49344
$ws.Range("L118").Value = 49344
$ws.Range("N118").Value = -52658
$ws.Range("H119").Value = 54890
$ws.Range("J119").Value = 54890
$ws.Range("L119").Value = 54890

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H95").Value = 44624
$ws.Range("J95").Value = 44624
$ws.Range("L95").Value = 44624
$ws.Range("N95").Value = -50116
$ws.Range("H117").Value = 0
$ws.Range("J117").Value = 0
$ws.Range("N117").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H70").Value = 31972.166
$ws.Range("J70").Value = 31972.166
$ws.Range("L70").Value = 31972.166
$ws.Range("H73").Value = 31972.166
$ws.Range("J73").Value = 31972.166
$ws.Range("L73").Value = 31972.166
$ws.Range("H106").Value = 35581.5
$ws.Range("J106").Value = 35581.5
$ws.Range("L106").Value = 35581.5
$ws.Range("N106").Value = -38105.5
$ws.Range("H111").Value = 47210
$ws.Range("J111").Value = 47210
$ws.Range("L111").Value = 47210
$ws.Range("N111").Value = -55390
$ws.Range("H116").Value = 47659.5
$ws.Range("J116").Value = 47659.5
$ws.Range("L116").Value = 47659.5
$ws.Range("N116").Value = -56837.5
$ws.Range("H125").Value = 49318
$ws.Range("J125").Value = 49318
$ws.Range("L125").Value = 49318
$ws.Range("N125").Value = -54238

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 744.4666999999999
$ws.Range("J4").Value = 1174.1111
$ws.Range("L4").Value = 3522.3333
$ws.Range("N4").Value = -3746.3333
$ws.Range("H25").Value = 900
$ws.Range("I25").Value = 800
$ws.Range("K25").Value = 2400
$ws.Range("H30").Value = 900
$ws.Range("I30").Value = 800
$ws.Range("K30").Value = 2400
$ws.Range("H46").Value = 1680.8
$ws.Range("J46").Value = 1680.8
$ws.Range("L46").Value = 5042.4
$ws.Range("N46").Value = -5224.4
$ws.Range("H56").Value = 8641.546
$ws.Range("I56").Value = 8641.546
$ws.Range("K56").Value = 8641.546
$ws.Range("M56").Value = -8111.546
$ws.Range("H113").Value = 5256.409
$ws.Range("I113").Value = 14978.857
$ws.Range("J113").Value = 719.26666
$ws.Range("K113").Value = 44936.571
$ws.Range("L113").Value = 2157.79998
$ws.Range("M113").Value = -42766.571
$ws.Range("N113").Value = -6497.79998
$ws.Range("H131").Value = 3255.2979
$ws.Range("I131").Value = 7223.7334
$ws.Range("J131").Value = 1395.0938
$ws.Range("K131").Value = 21671.2002
$ws.Range("L131").Value = 4185.2814
$ws.Range("M131").Value = -16631.2002
$ws.Range("N131").Value = -14265.2814

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H104").Value = 44878.25
$ws.Range("J104").Value = 44878.25
$ws.Range("L104").Value = 44878.25
$ws.Range("N104").Value = -51866.25
$ws.Range("H105").Value = 42933
$ws.Range("J105").Value = 42933
$ws.Range("L105").Value = 42933
$ws.Range("N105").Value = -49921
$ws.Range("H116").Value = 38912.285
$ws.Range("J116").Value = 38912.285
$ws.Range("L116").Value = 38912.285
$ws.Range("N116").Value = -48090.285
$ws.Range("H118").Value = 38183.332
$ws.Range("J118").Value = 38183.332
$ws.Range("L118").Value = 38183.332
$ws.Range("N118").Value = -41497.332
$ws.Range("H120").Value = 28766.666
$ws.Range("J120").Value = 28766.666
$ws.Range("L120").Value = 28766.666
$ws.Range("N120").Value = -38442.666
$ws.Range("H130").Value = 44523.5
$ws.Range("J130").Value = 44523.5
$ws.Range("L130").Value = 44523.5
$ws.Range("N130").Value = -54563.5
$ws.Range("H131").Value = 40658
$ws.Range("J131").Value = 40658
$ws.Range("L131").Value = 40658
$ws.Range("N131").Value = -50738

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H110").Value = 45499
$ws.Range("J110").Value = 45499
$ws.Range("L110").Value = 45499
$ws.Range("N110").Value = -53679

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H18").Value = 9761.666999999999
$ws.Range("I18").Value = 8570
$ws.Range("K18").Value = 8570
$ws.Range("M18").Value = -8397
$ws.Range("H27").Value = 36992
$ws.Range("J27").Value = 36992
$ws.Range("L27").Value = 36992
$ws.Range("N27").Value = -37130
$ws.Range("H103").Value = 42971
$ws.Range("J103").Value = 42971
$ws.Range("L103").Value = 42971
$ws.Range("N103").Value = -45315
$ws.Range("H105").Value = 50045
$ws.Range("J105").Value = 50045
$ws.Range("L105").Value = 50045
$ws.Range("N105").Value = -57033
$ws.Range("H109").Value = 23626
$ws.Range("J109").Value = 23626
$ws.Range("L109").Value = 23626
$ws.Range("N109").Value = -26400
$ws.Range("H129").Value = 39421
$ws.Range("J129").Value = 39421
$ws.Range("L129").Value = 39421
$ws.Range("N129").Value = -49421
